$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team-specific transition matrix values (Marshall_A)
$ws.Range("B2").Value = 0.2218649517684887
$ws.Range("C2").Value = 0.5112540192926045
$ws.Range("J2").Value = 0.009646302250803859
$ws.Range("P2").Value = 0.1543408360128617
$ws.Range("S2").Value = 0.1028938906752412
$ws.Range("B3").Value = 0.01204819277108434
$ws.Range("C3").Value = 0.03012048192771084
$ws.Range("J3").Value = 0.01807228915662651
$ws.Range("P3").Value = 0.7228915662650602
$ws.Range("S3").Value = 0.2168674698795181
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.3157894736842105
$ws.Range("B6").Value = 0.02512562814070352
$ws.Range("D6").Value = 0.02010050251256281
$ws.Range("F6").Value = 0.05025125628140704
$ws.Range("J6").Value = 0.2814070351758794
$ws.Range("O6").Value = 0.01507537688442211
$ws.Range("Q6").Value = 0.1557788944723618
$ws.Range("R6").Value = 0.08040201005025126
$ws.Range("S6").Value = 0.3718592964824121
$ws.Range("B7").Value = 0.1256830601092896
$ws.Range("D7").Value = 0.00546448087431694
$ws.Range("F7").Value = 0.0273224043715847
$ws.Range("J7").Value = 0.1639344262295082
$ws.Range("Q7").Value = 0.185792349726776
$ws.Range("R7").Value = 0.08743169398907104
$ws.Range("S7").Value = 0.4043715846994536
$ws.Range("B8").Value = 0.07436399217221135
$ws.Range("D8").Value = 0.0136986301369863
$ws.Range("E8").Value = 0.001956947162426614
$ws.Range("F8").Value = 0.04892367906066536
$ws.Range("J8").Value = 0.1389432485322896
$ws.Range("O8").Value = 0.01956947162426614
$ws.Range("Q8").Value = 0.1565557729941291
$ws.Range("R8").Value = 0.1095890410958904
$ws.Range("S8").Value = 0.436399217221135
$ws.Range("B9").Value = 0.07547169811320754
$ws.Range("D9").Value = 0.01415094339622642
$ws.Range("F9").Value = 0.05188679245283019
$ws.Range("J9").Value = 0.1179245283018868
$ws.Range("O9").Value = 0.02358490566037736
$ws.Range("Q9").Value = 0.1886792452830189
$ws.Range("R9").Value = 0.09905660377358491
$ws.Range("S9").Value = 0.4292452830188679
$ws.Range("B10").Value = 0.1111913357400722
$ws.Range("D10").Value = 0.01732851985559567
$ws.Range("E10").Value = 0.001444043321299639
$ws.Range("F10").Value = 0.06642599277978339
$ws.Range("J10").Value = 0.1299638989169675
$ws.Range("O10").Value = 0.005054151624548736
$ws.Range("Q10").Value = 0.1884476534296029
$ws.Range("R10").Value = 0.08014440433212996
$ws.Range("S10").Value = 0.4
$ws.Range("G11").Value = 0.1122448979591837
$ws.Range("J11").Value = 0.09523809523809523
$ws.Range("K11").Value = 0.1530612244897959
$ws.Range("L11").Value = 0.6360544217687075
$ws.Range("S11").Value = 0.003401360544217687
$ws.Range("G12").Value = 0.7058823529411765
$ws.Range("J12").Value = 0.267379679144385
$ws.Range("S12").Value = 0.0267379679144385
$ws.Range("G13").Value = 0.6111111111111112
$ws.Range("J13").Value = 0.3611111111111111
$ws.Range("S13").Value = 0.02777777777777778
$ws.Range("F15").Value = 0.01990049751243781
$ws.Range("H15").Value = 0.1791044776119403
$ws.Range("I15").Value = 0.09950248756218906
$ws.Range("J15").Value = 0.3432835820895522
$ws.Range("K15").Value = 0.0945273631840796
$ws.Range("M15").Value = 0.009950248756218905
$ws.Range("O15").Value = 0.03482587064676617
$ws.Range("S15").Value = 0.2189054726368159
$ws.Range("F16").Value = 0.02072538860103627
$ws.Range("H16").Value = 0.2176165803108808
$ws.Range("I16").Value = 0.06217616580310881
$ws.Range("J16").Value = 0.4455958549222798
$ws.Range("K16").Value = 0.06735751295336788
$ws.Range("M16").Value = 0.0155440414507772
$ws.Range("O16").Value = 0.04663212435233161
$ws.Range("S16").Value = 0.1243523316062176
$ws.Range("F17").Value = 0.006864988558352402
$ws.Range("H17").Value = 0.1807780320366133
$ws.Range("I17").Value = 0.1167048054919908
$ws.Range("J17").Value = 0.4073226544622426
$ws.Range("K17").Value = 0.08924485125858124
$ws.Range("M17").Value = 0.02288329519450801
$ws.Range("O17").Value = 0.04805491990846682
$ws.Range("S17").Value = 0.1281464530892449
$ws.Range("F18").Value = 0.01376146788990826
$ws.Range("H18").Value = 0.1834862385321101
$ws.Range("I18").Value = 0.05963302752293578
$ws.Range("J18").Value = 0.463302752293578
$ws.Range("K18").Value = 0.07339449541284404
$ws.Range("M18").Value = 0.009174311926605505
$ws.Range("O18").Value = 0.06422018348623854
$ws.Range("S18").Value = 0.1330275229357798
$ws.Range("F19").Value = 0.01375814627081825
$ws.Range("H19").Value = 0.2295438088341781
$ws.Range("I19").Value = 0.08472121650977553
$ws.Range("J19").Value = 0.3678493845039826
$ws.Range("K19").Value = 0.1151339608979001
$ws.Range("M19").Value = 0.01448225923244026
$ws.Range("N19").Value = 0.002172338884866039
$ws.Range("O19").Value = 0.06879073135409124
$ws.Range("S19").Value = 0.1035481535119479
